# Generate Report for Handback
# - Mark the two "Ready for handoff" status cells as "Handed back: in sync with en-US"
# - Fill in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   columns (I/J/K) for the zh-cn and de-de localization tables, now that the files
#   have been handed back.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status column text ---------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 30
$wsOverview.Columns.Item(6).ColumnWidth = 30

# --- zh-cn sheet ----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "340b87f2-7574-4a56-a057-5b123b33059a.md"
$wsZh.Range("J2").Value = "340b87f2-7574-4a56-a057-5b123b33059a.46bc87270154c6dd4b7abbc0098aa818c9a7eac8.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-15 20:24:59"

$wsZh.Range("I3").Value = "52c24b06-4060-4e58-be25-8892481a16e9.md"
$wsZh.Range("J3").Value = "52c24b06-4060-4e58-be25-8892481a16e9.18ee8285afde44d2759b01b20b6370e9f4f1393a.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-15 20:24:59"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cfcf12831813b47ed43868ea61989f6d2e7b31e/e2e/340b87f2-7574-4a56-a057-5b123b33059a.md", "", "", "340b87f2-7574-4a56-a057-5b123b33059a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cfcf12831813b47ed43868ea61989f6d2e7b31e/e2e/52c24b06-4060-4e58-be25-8892481a16e9.md", "", "", "52c24b06-4060-4e58-be25-8892481a16e9.md")
$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Range("I3").Style = "Hyperlink"

$wsZh.Columns.Item(3).ColumnWidth = 30
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "340b87f2-7574-4a56-a057-5b123b33059a.md"
$wsDe.Range("J2").Value = "340b87f2-7574-4a56-a057-5b123b33059a.46bc87270154c6dd4b7abbc0098aa818c9a7eac8.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-15 20:25:16"

$wsDe.Range("G3").Value = "52c24b06-4060-4e58-be25-8892481a16e9.18ee8285afde44d2759b01b20b6370e9f4f1393a.de-de.xlf"
$wsDe.Range("I3").Value = "52c24b06-4060-4e58-be25-8892481a16e9.md"
$wsDe.Range("J3").Value = "52c24b06-4060-4e58-be25-8892481a16e9.18ee8285afde44d2759b01b20b6370e9f4f1393a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-15 20:25:16"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cfcf12831813b47ed43868ea61989f6d2e7b31e/e2e/340b87f2-7574-4a56-a057-5b123b33059a.md", "", "", "340b87f2-7574-4a56-a057-5b123b33059a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cfcf12831813b47ed43868ea61989f6d2e7b31e/e2e/52c24b06-4060-4e58-be25-8892481a16e9.md", "", "", "52c24b06-4060-4e58-be25-8892481a16e9.md")
$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Range("I3").Style = "Hyperlink"

$wsDe.Columns.Item(3).ColumnWidth = 30
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
